$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2)
$ws.Range("B2").Value = 385.58204211234107
$ws.Range("C2").Value = 459.71363799849394
$ws.Range("D2").Value = 382.28634882838139
$ws.Range("E2").Value = 459.2903454735561

# Update STR row (row 3)
$ws.Range("B3").Value = 388.3339396377682
$ws.Range("C3").Value = 471.86629880563208
$ws.Range("D3").Value = 388.74886006185744
$ws.Range("E3").Value = 463.22608326593712

# Update selection to match new selected range
$ws.Range("B1:E3").Select()
